# Mark the name "p2ysb4km" as used:
#  - remove it from the pool of available names on "Sheet1" (it currently
#    sits in row 2; deleting the row shifts everything below up by one,
#    shrinking the sheet from A1:A464 to A1:A463)
#  - append a record for it to the "used" sheet (new row 36), capturing the
#    source filename and the timestamp it was consumed at

$wb = $excel.ActiveWorkbook
$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet  = $wb.Worksheets.Item("used")

# Remove the consumed name from the names pool.
$namesSheet.Rows.Item(2).Delete()

# Record it as used on the "used" sheet.
$nextRow = $usedSheet.Cells.Item($usedSheet.Rows.Count, 1).End(-4162).Row + 1
if ($nextRow -lt 2) { $nextRow = 2 }

$usedSheet.Cells.Item($nextRow, 1).Value = "p2ysb4km"
$usedSheet.Cells.Item($nextRow, 2).Value = "ChatGPT Image 2026年1月21日 16_39_42.png"
$usedSheet.Cells.Item($nextRow, 3).Value = "2026-01-21 16:40:40"
